$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on price cells whose new values would otherwise
# be auto-detected as numbers by Excel (losing the original text formatting).
$textCells = @("D5", "D6", "D8", "D10", "D11", "D13", "D14", "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D29", "D30", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "63.228.70"
$ws.Range("E2").Value = "  +6.32%  "

$ws.Range("D3").Value = "2.431.81"
$ws.Range("E3").Value = "  +6.02%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "566.74"
$ws.Range("E5").Value = "  +4.82%  "

$ws.Range("D6").Value = "140.95"
$ws.Range("E6").Value = "  +9.57%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "0.586"
$ws.Range("E8").Value = "  +2.99%  "

$ws.Range("D9").Value = "2.426.82"
$ws.Range("E9").Value = "  +5.87%  "

$ws.Range("D10").Value = "0.105"
$ws.Range("E10").Value = "  +4.71%  "

$ws.Range("D11").Value = "5.74"
$ws.Range("E11").Value = "  +4.02%  "

$ws.Range("E12").Value = "  +0.56%  "

$ws.Range("D13").Value = "0.351"
$ws.Range("E13").Value = "  +6.52%  "

$ws.Range("D14").Value = "26.32"
$ws.Range("E14").Value = "  +14.13%  "

$ws.Range("D15").Value = "2.866.41"
$ws.Range("E15").Value = "  +6.04%  "

$ws.Range("D16").Value = "63.038.96"
$ws.Range("E16").Value = "  +6.06%  "

$ws.Range("D17").Value = "0.0000143"
$ws.Range("E17").Value = "  +8.98%  "

$ws.Range("D18").Value = "2.427.72"
$ws.Range("E18").Value = "  +5.85%  "

$ws.Range("D19").Value = "11.24"
$ws.Range("E19").Value = "  +8.19%  "

$ws.Range("D20").Value = "340.18"
$ws.Range("E20").Value = "  +10.09%  "

$ws.Range("D21").Value = "4.23"
$ws.Range("E21").Value = "  +5.18%  "

$ws.Range("D22").Value = "6.80"
$ws.Range("E22").Value = "  +4.79%  "

$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.02%  "

$ws.Range("D24").Value = "5.65"
$ws.Range("E24").Value = "  +0.30%  "

$ws.Range("D25").Value = "65.38"
$ws.Range("E25").Value = "  +4.30%  "

$ws.Range("D26").Value = "0.175"
$ws.Range("E26").Value = "  +4.37%  "

$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.16%  "

$ws.Range("E28").Value = "  +14.47%  "

$ws.Range("D29").Value = "8.20"
$ws.Range("E29").Value = "  +6.75%  "

$ws.Range("D30").Value = "1.35"
$ws.Range("E30").Value = "  +13.42%  "

$ws.Range("D31").Value = "0.0₃0797"
$ws.Range("E31").Value = "  +11.44%  "

$ws.Range("E32").Value = "  +7.05%  "

$ws.Range("D33").Value = "6.56"
$ws.Range("E33").Value = "  +13.83%  "

$ws.Range("D34").Value = "174.12"
$ws.Range("E34").Value = "  +1.33%  "

$ws.Range("D35").Value = "1.49"
$ws.Range("E35").Value = "  +12.09%  "

$ws.Range("D36").Value = "0.398"
$ws.Range("E36").Value = "  +5.74%  "

$ws.Range("D37").Value = "18.63"
$ws.Range("E37").Value = "  +6.21%  "

$ws.Range("D38").Value = "371.19"
$ws.Range("E38").Value = "  +19.29%  "

$ws.Range("D39").Value = "4.48"
$ws.Range("E39").Value = "  +12.82%  "

$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.03%  "

$ws.Range("D42").Value = "1.70"
$ws.Range("E42").Value = "  +14.42%  "

$ws.Range("D43").Value = "39.96"
$ws.Range("E43").Value = "  +6.93%  "

$ws.Range("D44").Value = "146.74"
$ws.Range("E44").Value = "  +7.83%  "

$ws.Range("D45").Value = "3.69"
$ws.Range("E45").Value = "  +8.35%  "

$ws.Range("D46").Value = "20.49"
$ws.Range("E46").Value = "  +11.15%  "

$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "0.0957"
$ws.Range("E47").Value = "  +2.02%  "

$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "0.592"
$ws.Range("E48").Value = "  +5.45%  "

$ws.Range("D49").Value = "0.0521"
$ws.Range("E49").Value = "  +7.20%  "

$ws.Range("D50").Value = "0.0225"
$ws.Range("E50").Value = "  +6.26%  "

$ws.Range("D51").Value = "17.86"
$ws.Range("E51").Value = "  +7.46%  "
